$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-03-23 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-24 Friday", 2) | Out-Null
$d.Content.Find.Execute("42×72=3024", $true, $false, $false, $false, $false, $true, 1, $false, "40×26=1040", 2) | Out-Null
$d.Content.Find.Execute("24×88=2112", $true, $false, $false, $false, $false, $true, 1, $false, "72×14=1008", 2) | Out-Null
$d.Content.Find.Execute("30×26=780", $true, $false, $false, $false, $false, $true, 1, $false, "73×43=3139", 2) | Out-Null
$d.Content.Find.Execute("81×23=1863", $true, $false, $false, $false, $false, $true, 1, $false, "60×51=3060", 2) | Out-Null
$d.Content.Find.Execute("93×22=2046", $true, $false, $false, $false, $false, $true, 1, $false, "82×39=3198", 2) | Out-Null
$d.Content.Find.Execute("43×29=1247", $true, $false, $false, $false, $false, $true, 1, $false, "16×32=512", 2) | Out-Null
$d.Content.Find.Execute("77×94=7238", $true, $false, $false, $false, $false, $true, 1, $false, "88×78=6864", 2) | Out-Null
$d.Content.Find.Execute("41×74=3034", $true, $false, $false, $false, $false, $true, 1, $false, "20×43=860", 2) | Out-Null
$d.Content.Find.Execute("90×11=990", $true, $false, $false, $false, $false, $true, 1, $false, "34×48=1632", 2) | Out-Null
$d.Content.Find.Execute("25×10=250", $true, $false, $false, $false, $false, $true, 1, $false, "84×20=1680", 2) | Out-Null
$d.Content.Find.Execute("45×99=4455", $true, $false, $false, $false, $false, $true, 1, $false, "63×70=4410", 2) | Out-Null
$d.Content.Find.Execute("53×40=2120", $true, $false, $false, $false, $false, $true, 1, $false, "44×22=968", 2) | Out-Null
$d.Content.Find.Execute("58×54=3132", $true, $false, $false, $false, $false, $true, 1, $false, "20×29=580", 2) | Out-Null
$d.Content.Find.Execute("97×38=3686", $true, $false, $false, $false, $false, $true, 1, $false, "84×22=1848", 2) | Out-Null
$d.Content.Find.Execute("41×36=1476", $true, $false, $false, $false, $false, $true, 1, $false, "47×88=4136", 2) | Out-Null
$d.Content.Find.Execute("15×54=810", $true, $false, $false, $false, $false, $true, 1, $false, "42×65=2730", 2) | Out-Null
$d.Content.Find.Execute("15×85=1275", $true, $false, $false, $false, $false, $true, 1, $false, "83×55=4565", 2) | Out-Null
$d.Content.Find.Execute("62×20=1240", $true, $false, $false, $false, $false, $true, 1, $false, "49×73=3577", 2) | Out-Null
$d.Content.Find.Execute("84×51=4284", $true, $false, $false, $false, $false, $true, 1, $false, "20×68=1360", 2) | Out-Null
$d.Content.Find.Execute("11×59=649", $true, $false, $false, $false, $false, $true, 1, $false, "21×29=609", 2) | Out-Null
$d.Content.Find.Execute("53×80=4240", $true, $false, $false, $false, $false, $true, 1, $false, "56×72=4032", 2) | Out-Null
$d.Content.Find.Execute("45×84=3780", $true, $false, $false, $false, $false, $true, 1, $false, "25×24=600", 2) | Out-Null
$d.Content.Find.Execute("34×18=612", $true, $false, $false, $false, $false, $true, 1, $false, "75×99=7425", 2) | Out-Null
$d.Content.Find.Execute("29×94=2726", $true, $false, $false, $false, $false, $true, 1, $false, "42×68=2856", 2) | Out-Null
$d.Content.Find.Execute("83×13=1079", $true, $false, $false, $false, $false, $true, 1, $false, "74×86=6364", 2) | Out-Null
$d.Content.Find.Execute("52×33=1716", $true, $false, $false, $false, $false, $true, 1, $false, "29×25=725", 2) | Out-Null
$d.Content.Find.Execute("76×79=6004", $true, $false, $false, $false, $false, $true, 1, $false, "97×31=3007", 2) | Out-Null
$d.Content.Find.Execute("59×41=2419", $true, $false, $false, $false, $false, $true, 1, $false, "94×48=4512", 2) | Out-Null
$d.Content.Find.Execute("77×20=1540", $true, $false, $false, $false, $false, $true, 1, $false, "82×44=3608", 2) | Out-Null
$d.Content.Find.Execute("79×46=3634", $true, $false, $false, $false, $false, $true, 1, $false, "91×28=2548", 2) | Out-Null
$d.Content.Find.Execute("57×40=2280", $true, $false, $false, $false, $false, $true, 1, $false, "11×41=451", 2) | Out-Null
$d.Content.Find.Execute("86×70=6020", $true, $false, $false, $false, $false, $true, 1, $false, "53×98=5194", 2) | Out-Null
$d.Content.Find.Execute("64×33=2112", $true, $false, $false, $false, $false, $true, 1, $false, "74×11=814", 2) | Out-Null
$d.Content.Find.Execute("16×75=1200", $true, $false, $false, $false, $false, $true, 1, $false, "75×99=7425", 2) | Out-Null
$d.Content.Find.Execute("96×34=3264", $true, $false, $false, $false, $false, $true, 1, $false, "30×90=2700", 2) | Out-Null
$d.Content.Find.Execute("36×90=3240", $true, $false, $false, $false, $false, $true, 1, $false, "37×17=629", 2) | Out-Null
$d.Content.Find.Execute("34×59=2006", $true, $false, $false, $false, $false, $true, 1, $false, "68×93=6324", 2) | Out-Null
$d.Content.Find.Execute("16×55=880", $true, $false, $false, $false, $false, $true, 1, $false, "95×85=8075", 2) | Out-Null
$d.Content.Find.Execute("24×70=1680", $true, $false, $false, $false, $false, $true, 1, $false, "38×38=1444", 2) | Out-Null
$d.Content.Find.Execute("86×69=5934", $true, $false, $false, $false, $false, $true, 1, $false, "10×84=840", 2) | Out-Null
$d.Content.Find.Execute("92×14=1288", $true, $false, $false, $false, $false, $true, 1, $false, "72×94=6768", 2) | Out-Null
$d.Content.Find.Execute("54×50=2700", $true, $false, $false, $false, $false, $true, 1, $false, "87×31=2697", 2) | Out-Null
$d.Content.Find.Execute("32×100=3200", $true, $false, $false, $false, $false, $true, 1, $false, "63×36=2268", 2) | Out-Null
$d.Content.Find.Execute("30×99=2970", $true, $false, $false, $false, $false, $true, 1, $false, "45×18=810", 2) | Out-Null
$d.Content.Find.Execute("73×24=1752", $true, $false, $false, $false, $false, $true, 1, $false, "82×97=7954", 2) | Out-Null
$d.Content.Find.Execute("54×61=3294", $true, $false, $false, $false, $false, $true, 1, $false, "72×22=1584", 2) | Out-Null
$d.Content.Find.Execute("26×100=2600", $true, $false, $false, $false, $false, $true, 1, $false, "22×80=1760", 2) | Out-Null
$d.Content.Find.Execute("89×40=3560", $true, $false, $false, $false, $false, $true, 1, $false, "99×84=8316", 2) | Out-Null
$d.Content.Find.Execute("15×34=510", $true, $false, $false, $false, $false, $true, 1, $false, "92×28=2576", 2) | Out-Null
$d.Content.Find.Execute("41×78=3198", $true, $false, $false, $false, $false, $true, 1, $false, "28×74=2072", 2) | Out-Null
$d.Content.Find.Execute("39×92=3588", $true, $false, $false, $false, $false, $true, 1, $false, "97×80=7760", 2) | Out-Null
$d.Content.Find.Execute("39×42=1638", $true, $false, $false, $false, $false, $true, 1, $false, "18×67=1206", 2) | Out-Null
$d.Content.Find.Execute("77×24=1848", $true, $false, $false, $false, $false, $true, 1, $false, "98×66=6468", 2) | Out-Null
$d.Content.Find.Execute("83×60=4980", $true, $false, $false, $false, $false, $true, 1, $false, "16×32=512", 2) | Out-Null
$d.Content.Find.Execute("56×97=5432", $true, $false, $false, $false, $false, $true, 1, $false, "73×25=1825", 2) | Out-Null
$d.Content.Find.Execute("37×26=962", $true, $false, $false, $false, $false, $true, 1, $false, "18×66=1188", 2) | Out-Null
$d.Content.Find.Execute("59×79=4661", $true, $false, $false, $false, $false, $true, 1, $false, "37×42=1554", 2) | Out-Null
$d.Content.Find.Execute("46×39=1794", $true, $false, $false, $false, $false, $true, 1, $false, "72×31=2232", 2) | Out-Null
$d.Content.Find.Execute("47×86=4042", $true, $false, $false, $false, $false, $true, 1, $false, "74×10=740", 2) | Out-Null
$d.Content.Find.Execute("23×57=1311", $true, $false, $false, $false, $false, $true, 1, $false, "78×99=7722", 2) | Out-Null
$d.Content.Find.Execute("12×52=624", $true, $false, $false, $false, $false, $true, 1, $false, "57×22=1254", 2) | Out-Null
$d.Content.Find.Execute("51×42=2142", $true, $false, $false, $false, $false, $true, 1, $false, "81×92=7452", 2) | Out-Null
$d.Content.Find.Execute("30×33=990", $true, $false, $false, $false, $false, $true, 1, $false, "51×88=4488", 2) | Out-Null
$d.Content.Find.Execute("52×10=520", $true, $false, $false, $false, $false, $true, 1, $false, "32×78=2496", 2) | Out-Null
$d.Content.Find.Execute("68×42=2856", $true, $false, $false, $false, $false, $true, 1, $false, "18×48=864", 2) | Out-Null
$d.Content.Find.Execute("65×10=650", $true, $false, $false, $false, $false, $true, 1, $false, "74×46=3404", 2) | Out-Null
$d.Content.Find.Execute("43×24=1032", $true, $false, $false, $false, $false, $true, 1, $false, "37×97=3589", 2) | Out-Null
$d.Content.Find.Execute("14×30=420", $true, $false, $false, $false, $false, $true, 1, $false, "91×56=5096", 2) | Out-Null
$d.Content.Find.Execute("73×86=6278", $true, $false, $false, $false, $false, $true, 1, $false, "26×31=806", 2) | Out-Null
$d.Content.Find.Execute("97×70=6790", $true, $false, $false, $false, $false, $true, 1, $false, "42×10=420", 2) | Out-Null
$d.Content.Find.Execute("88×72=6336", $true, $false, $false, $false, $false, $true, 1, $false, "69×79=5451", 2) | Out-Null
$d.Content.Find.Execute("33×97=3201", $true, $false, $false, $false, $false, $true, 1, $false, "30×62=1860", 2) | Out-Null
$d.Content.Find.Execute("73×23=1679", $true, $false, $false, $false, $false, $true, 1, $false, "27×37=999", 2) | Out-Null
$d.Content.Find.Execute("66×11=726", $true, $false, $false, $false, $false, $true, 1, $false, "67×81=5427", 2) | Out-Null
$d.Content.Find.Execute("12×91=1092", $true, $false, $false, $false, $false, $true, 1, $false, "49×100=4900", 2) | Out-Null
$d.Content.Find.Execute("38×21=798", $true, $false, $false, $false, $false, $true, 1, $false, "77×84=6468", 2) | Out-Null
$d.Content.Find.Execute("80×39=3120", $true, $false, $false, $false, $false, $true, 1, $false, "40×11=440", 2) | Out-Null
$d.Content.Find.Execute("51×36=1836", $true, $false, $false, $false, $false, $true, 1, $false, "59×65=3835", 2) | Out-Null
$d.Content.Find.Execute("82×45=3690", $true, $false, $false, $false, $false, $true, 1, $false, "40×63=2520", 2) | Out-Null
$d.Content.Find.Execute("40×100=4000", $true, $false, $false, $false, $false, $true, 1, $false, "62×12=744", 2) | Out-Null
$d.Content.Find.Execute("92×70=6440", $true, $false, $false, $false, $false, $true, 1, $false, "25×89=2225", 2) | Out-Null
$d.Content.Find.Execute("53×67=3551", $true, $false, $false, $false, $false, $true, 1, $false, "32×25=800", 2) | Out-Null
$d.Content.Find.Execute("50×12=600", $true, $false, $false, $false, $false, $true, 1, $false, "47×11=517", 2) | Out-Null
$d.Content.Find.Execute("75×82=6150", $true, $false, $false, $false, $false, $true, 1, $false, "63×32=2016", 2) | Out-Null
$d.Content.Find.Execute("70×97=6790", $true, $false, $false, $false, $false, $true, 1, $false, "72×36=2592", 2) | Out-Null
$d.Content.Find.Execute("67×61=4087", $true, $false, $false, $false, $false, $true, 1, $false, "75×38=2850", 2) | Out-Null
$d.Content.Find.Execute("11×75=825", $true, $false, $false, $false, $false, $true, 1, $false, "56×35=1960", 2) | Out-Null
$d.Content.Find.Execute("98×85=8330", $true, $false, $false, $false, $false, $true, 1, $false, "90×64=5760", 2) | Out-Null
$d.Content.Find.Execute("58×91=5278", $true, $false, $false, $false, $false, $true, 1, $false, "88×17=1496", 2) | Out-Null
$d.Content.Find.Execute("38×50=1900", $true, $false, $false, $false, $false, $true, 1, $false, "99×92=9108", 2) | Out-Null
$d.Content.Find.Execute("21×83=1743", $true, $false, $false, $false, $false, $true, 1, $false, "61×64=3904", 2) | Out-Null
$d.Content.Find.Execute("73×80=5840", $true, $false, $false, $false, $false, $true, 1, $false, "69×66=4554", 2) | Out-Null
$d.Content.Find.Execute("45×97=4365", $true, $false, $false, $false, $false, $true, 1, $false, "87×48=4176", 2) | Out-Null
$d.Content.Find.Execute("80×89=7120", $true, $false, $false, $false, $false, $true, 1, $false, "20×78=1560", 2) | Out-Null
$d.Content.Find.Execute("20×30=600", $true, $false, $false, $false, $false, $true, 1, $false, "99×90=8910", 2) | Out-Null
$d.Content.Find.Execute("91×31=2821", $true, $false, $false, $false, $false, $true, 1, $false, "83×10=830", 2) | Out-Null
$d.Content.Find.Execute("66×18=1188", $true, $false, $false, $false, $false, $true, 1, $false, "66×79=5214", 2) | Out-Null
$d.Content.Find.Execute("86×85=7310", $true, $false, $false, $false, $false, $true, 1, $false, "62×26=1612", 2) | Out-Null
$d.Content.Find.Execute("25×93=2325", $true, $false, $false, $false, $false, $true, 1, $false, "10×94=940", 2) | Out-Null
$d.Content.Find.Execute("42×15=630", $true, $false, $false, $false, $false, $true, 1, $false, "86×55=4730", 2) | Out-Null
